# Update cryptocurrency price/volume figures per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.600.01"
$ws.Range("E2").Value = "  -5.25%  "

$ws.Range("D3").Value = "2.207.45"
$ws.Range("E3").Value = "  -7.40%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  -4.89%  "

$ws.Range("D9").Value = "2.233.39"
$ws.Range("E9").Value = "  -6.86%  "

$ws.Range("E10").Value = "  -6.95%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.85%  "

$ws.Range("D14").Value = "2.603.10"
$ws.Range("E14").Value = "  -7.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.11%  "

$ws.Range("D16").Value = "53.533.73"
$ws.Range("E16").Value = "  -5.29%  "

$ws.Range("E17").Value = "  -4.44%  "

$ws.Range("D18").Value = "2.192.42"
$ws.Range("E18").Value = "  -8.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "295.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.86%  "

$ws.Range("E22").Value = "  -2.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.366"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.41%  "

$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").Value = "2.310.15"
$ws.Range("E28").Value = "  -7.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "162.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "0.0₃0670"
$ws.Range("E33").Value = "  -6.74%  "

$ws.Range("E34").Value = "  -1.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.43%  "

$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.367"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.25%  "

$ws.Range("E43").Value = "  -2.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "126.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.87%  "

$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0882"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.535"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "232.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.81%  "

$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0201"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
